$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"47.25342666666666"
$ws.Range("H2").Value = [double]"141.76028"
$ws.Range("I2").Value = [double]"0.7256581002375004"
$ws.Range("J2").Value = [double]"0.7256581002375005"
$ws.Range("M2").Value = [double]"6.305846"
$ws.Range("N2").Value = [double]"18.917538"
$ws.Range("O2").Value = [double]"0.01356150511917599"
$ws.Range("P2").Value = [double]"0.01356150511917599"
$ws.Range("Q2").Value = [double]"297.9728315322933"
$ws.Range("R2").Value = [double]"2681.75548379064"
$ws.Range("S2").Value = [double]"0.009841016041142384"
$ws.Range("T2").Value = [double]"0.009841016041142386"
$ws.Range("G3").Value = [double]"47.25342666666666"
$ws.Range("H3").Value = [double]"141.76028"
$ws.Range("I3").Value = [double]"0.7256581002375004"
$ws.Range("J3").Value = [double]"0.7256581002375005"
$ws.Range("O3").Value = [double]"0.392557056479861"
$ws.Range("P3").Value = [double]"0.3925570564798609"
$ws.Range("Q3").Value = [double]"8625.247465481461"
$ws.Range("R3").Value = [double]"77627.22718933316"
$ws.Range("S3").Value = [double]"0.2848622078400011"
$ws.Range("T3").Value = [double]"0.2848622078400011"
$ws.Range("G4").Value = [double]"47.25342666666666"
$ws.Range("H4").Value = [double]"141.76028"
$ws.Range("I4").Value = [double]"0.7256581002375004"
$ws.Range("J4").Value = [double]"0.7256581002375005"
$ws.Range("M4").Value = [double]"127.396393"
$ws.Range("N4").Value = [double]"382.189179"
$ws.Range("O4").Value = [double]"0.2739817680029065"
$ws.Range("P4").Value = [double]"0.2739817680029065"
$ws.Range("Q4").Value = [double]"6019.916114223346"
$ws.Range("R4").Value = [double]"54179.24502801012"
$ws.Range("S4").Value = [double]"0.1988170892687007"
$ws.Range("T4").Value = [double]"0.1988170892687008"
$ws.Range("G5").Value = [double]"47.25342666666666"
$ws.Range("H5").Value = [double]"141.76028"
$ws.Range("I5").Value = [double]"0.7256581002375004"
$ws.Range("J5").Value = [double]"0.7256581002375005"
$ws.Range("M5").Value = [double]"19.42400133333333"
$ws.Range("N5").Value = [double]"58.272004"
$ws.Range("O5").Value = [double]"0.04177372766745037"
$ws.Range("P5").Value = [double]"0.04177372766745036"
$ws.Range("Q5").Value = [double]"917.8506225779021"
$ws.Range("R5").Value = [double]"8260.655603201119"
$ws.Range("S5").Value = [double]"0.03031344385900074"
$ws.Range("T5").Value = [double]"0.03031344385900074"
$ws.Range("G6").Value = [double]"47.25342666666666"
$ws.Range("H6").Value = [double]"141.76028"
$ws.Range("I6").Value = [double]"0.7256581002375004"
$ws.Range("J6").Value = [double]"0.7256581002375005"
$ws.Range("M6").Value = [double]"129.3233566666667"
$ws.Range("N6").Value = [double]"387.97007"
$ws.Range("O6").Value = [double]"0.2781259427306063"
$ws.Range("P6").Value = [double]"0.2781259427306062"
$ws.Range("Q6").Value = [double]"6110.971750535511"
$ws.Range("R6").Value = [double]"54998.7457548196"
$ws.Range("S6").Value = [double]"0.2018243432286556"
$ws.Range("T6").Value = [double]"0.2018243432286556"
$ws.Range("I7").Value = [double]"0.0008031160114981568"
$ws.Range("J7").Value = [double]"0.0008031160114981569"
$ws.Range("M7").Value = [double]"6.305846"
$ws.Range("N7").Value = [double]"18.917538"
$ws.Range("O7").Value = [double]"0.01356150511917599"
$ws.Range("P7").Value = [double]"0.01356150511917599"
$ws.Range("Q7").Value = [double]"0.3297789302106667"
$ws.Range("R7").Value = [double]"2.968010371896"
$ws.Range("S7").Value = [double]"1.089146190122446E-05"
$ws.Range("T7").Value = [double]"1.089146190122446E-05"
$ws.Range("I8").Value = [double]"0.0008031160114981568"
$ws.Range("J8").Value = [double]"0.0008031160114981569"
$ws.Range("O8").Value = [double]"0.392557056479861"
$ws.Range("P8").Value = [double]"0.3925570564798609"
$ws.Range("S8").Value = [double]"0.0003152688574855626"
$ws.Range("T8").Value = [double]"0.0003152688574855626"
$ws.Range("I9").Value = [double]"0.0008031160114981568"
$ws.Range("J9").Value = [double]"0.0008031160114981569"
$ws.Range("M9").Value = [double]"127.396393"
$ws.Range("N9").Value = [double]"382.189179"
$ws.Range("O9").Value = [double]"0.2739817680029065"
$ws.Range("P9").Value = [double]"0.2739817680029065"
$ws.Range("Q9").Value = [double]"6.662491630185333"
$ws.Range("R9").Value = [double]"59.962424671668"
$ws.Range("S9").Value = [double]"0.0002200391447417076"
$ws.Range("T9").Value = [double]"0.0002200391447417076"
$ws.Range("I10").Value = [double]"0.0008031160114981568"
$ws.Range("J10").Value = [double]"0.0008031160114981569"
$ws.Range("M10").Value = [double]"19.42400133333333"
$ws.Range("N10").Value = [double]"58.272004"
$ws.Range("O10").Value = [double]"0.04177372766745037"
$ws.Range("P10").Value = [double]"0.04177372766745036"
$ws.Range("Q10").Value = [double]"1.015823472396445"
$ws.Range("R10").Value = [double]"9.142411251568001"
$ws.Range("S10").Value = [double]"3.354914954969294E-05"
$ws.Range("T10").Value = [double]"3.354914954969294E-05"
$ws.Range("I11").Value = [double]"0.0008031160114981568"
$ws.Range("J11").Value = [double]"0.0008031160114981569"
$ws.Range("M11").Value = [double]"129.3233566666667"
$ws.Range("N11").Value = [double]"387.97007"
$ws.Range("O11").Value = [double]"0.2781259427306063"
$ws.Range("P11").Value = [double]"0.2781259427306062"
$ws.Range("Q11").Value = [double]"6.763266691382223"
$ws.Range("R11").Value = [double]"60.86940022244001"
$ws.Range("S11").Value = [double]"0.0002233673978199693"
$ws.Range("T11").Value = [double]"0.0002233673978199692"
$ws.Range("G12").Value = [double]"6.996562666666667"
$ws.Range("H12").Value = [double]"20.989688"
$ws.Range("I12").Value = [double]"0.1074443216298519"
$ws.Range("J12").Value = [double]"0.1074443216298519"
$ws.Range("M12").Value = [double]"6.305846"
$ws.Range("N12").Value = [double]"18.917538"
$ws.Range("O12").Value = [double]"0.01356150511917599"
$ws.Range("P12").Value = [double]"0.01356150511917599"
$ws.Range("Q12").Value = [double]"44.11924670534933"
$ws.Range("R12").Value = [double]"397.073220348144"
$ws.Range("S12").Value = [double]"0.001457106717809628"
$ws.Range("T12").Value = [double]"0.001457106717809628"
$ws.Range("G13").Value = [double]"6.996562666666667"
$ws.Range("H13").Value = [double]"20.989688"
$ws.Range("I13").Value = [double]"0.1074443216298519"
$ws.Range("J13").Value = [double]"0.1074443216298519"
$ws.Range("O13").Value = [double]"0.392557056479861"
$ws.Range("P13").Value = [double]"0.3925570564798609"
$ws.Range("Q13").Value = [double]"1277.09435409726"
$ws.Range("R13").Value = [double]"11493.84918687534"
$ws.Range("S13").Value = [double]"0.04217802663449012"
$ws.Range("T13").Value = [double]"0.04217802663449012"
$ws.Range("G14").Value = [double]"6.996562666666667"
$ws.Range("H14").Value = [double]"20.989688"
$ws.Range("I14").Value = [double]"0.1074443216298519"
$ws.Range("J14").Value = [double]"0.1074443216298519"
$ws.Range("M14").Value = [double]"127.396393"
$ws.Range("N14").Value = [double]"382.189179"
$ws.Range("O14").Value = [double]"0.2739817680029065"
$ws.Range("P14").Value = [double]"0.2739817680029065"
$ws.Range("Q14").Value = [double]"891.3368471317947"
$ws.Range("R14").Value = [double]"8022.031624186153"
$ws.Range("S14").Value = [double]"0.02943778520201976"
$ws.Range("T14").Value = [double]"0.02943778520201976"
$ws.Range("G15").Value = [double]"6.996562666666667"
$ws.Range("H15").Value = [double]"20.989688"
$ws.Range("I15").Value = [double]"0.1074443216298519"
$ws.Range("J15").Value = [double]"0.1074443216298519"
$ws.Range("M15").Value = [double]"19.42400133333333"
$ws.Range("N15").Value = [double]"58.272004"
$ws.Range("O15").Value = [double]"0.04177372766745037"
$ws.Range("P15").Value = [double]"0.04177372766745036"
$ws.Range("Q15").Value = [double]"135.9012425660836"
$ws.Range("R15").Value = [double]"1223.111183094752"
$ws.Range("S15").Value = [double]"0.004488349831179381"
$ws.Range("T15").Value = [double]"0.004488349831179381"
$ws.Range("G16").Value = [double]"6.996562666666667"
$ws.Range("H16").Value = [double]"20.989688"
$ws.Range("I16").Value = [double]"0.1074443216298519"
$ws.Range("J16").Value = [double]"0.1074443216298519"
$ws.Range("M16").Value = [double]"129.3233566666667"
$ws.Range("N16").Value = [double]"387.97007"
$ws.Range("O16").Value = [double]"0.2781259427306063"
$ws.Range("P16").Value = [double]"0.2781259427306062"
$ws.Range("Q16").Value = [double]"904.8189691820179"
$ws.Range("R16").Value = [double]"8143.370722638161"
$ws.Range("S16").Value = [double]"0.02988305324435303"
$ws.Range("T16").Value = [double]"0.02988305324435303"
$ws.Range("G17").Value = [double]"2.605199666666667"
$ws.Range("H17").Value = [double]"7.815599000000001"
$ws.Range("I17").Value = [double]"0.0400073470689964"
$ws.Range("J17").Value = [double]"0.0400073470689964"
$ws.Range("M17").Value = [double]"6.305846"
$ws.Range("N17").Value = [double]"18.917538"
$ws.Range("O17").Value = [double]"0.01356150511917599"
$ws.Range("P17").Value = [double]"0.01356150511917599"
$ws.Range("Q17").Value = [double]"16.42798789725133"
$ws.Range("R17").Value = [double]"147.851891075262"
$ws.Range("S17").Value = [double]"0.0005425598420808452"
$ws.Range("T17").Value = [double]"0.0005425598420808451"
$ws.Range("G18").Value = [double]"2.605199666666667"
$ws.Range("H18").Value = [double]"7.815599000000001"
$ws.Range("I18").Value = [double]"0.0400073470689964"
$ws.Range("J18").Value = [double]"0.0400073470689964"
$ws.Range("O18").Value = [double]"0.392557056479861"
$ws.Range("P18").Value = [double]"0.3925570564798609"
$ws.Range("Q18").Value = [double]"475.531477970906"
$ws.Range("R18").Value = [double]"4279.783301738154"
$ws.Range("S18").Value = [double]"0.01570516640297342"
$ws.Range("T18").Value = [double]"0.01570516640297342"
$ws.Range("G19").Value = [double]"2.605199666666667"
$ws.Range("H19").Value = [double]"7.815599000000001"
$ws.Range("I19").Value = [double]"0.0400073470689964"
$ws.Range("J19").Value = [double]"0.0400073470689964"
$ws.Range("M19").Value = [double]"127.396393"
$ws.Range("N19").Value = [double]"382.189179"
$ws.Range("O19").Value = [double]"0.2739817680029065"
$ws.Range("P19").Value = [double]"0.2739817680029065"
$ws.Range("Q19").Value = [double]"331.8930405781357"
$ws.Range("R19").Value = [double]"2987.037365203221"
$ws.Range("S19").Value = [double]"0.01096128368306954"
$ws.Range("T19").Value = [double]"0.01096128368306954"
$ws.Range("G20").Value = [double]"2.605199666666667"
$ws.Range("H20").Value = [double]"7.815599000000001"
$ws.Range("I20").Value = [double]"0.0400073470689964"
$ws.Range("J20").Value = [double]"0.0400073470689964"
$ws.Range("M20").Value = [double]"19.42400133333333"
$ws.Range("N20").Value = [double]"58.272004"
$ws.Range("O20").Value = [double]"0.04177372766745037"
$ws.Range("P20").Value = [double]"0.04177372766745036"
$ws.Range("Q20").Value = [double]"50.60340179893289"
$ws.Range("R20").Value = [double]"455.4306161903961"
$ws.Range("S20").Value = [double]"0.001671256021157425"
$ws.Range("T20").Value = [double]"0.001671256021157424"
$ws.Range("G21").Value = [double]"2.605199666666667"
$ws.Range("H21").Value = [double]"7.815599000000001"
$ws.Range("I21").Value = [double]"0.0400073470689964"
$ws.Range("J21").Value = [double]"0.0400073470689964"
$ws.Range("M21").Value = [double]"129.3233566666667"
$ws.Range("N21").Value = [double]"387.97007"
$ws.Range("O21").Value = [double]"0.2781259427306063"
$ws.Range("P21").Value = [double]"0.2781259427306062"
$ws.Range("Q21").Value = [double]"336.9131656802145"
$ws.Range("R21").Value = [double]"3032.218491121931"
$ws.Range("S21").Value = [double]"0.01112708111971518"
$ws.Range("T21").Value = [double]"0.01112708111971518"
$ws.Range("G22").Value = [double]"8.210544666666665"
$ws.Range("H22").Value = [double]"24.631634"
$ws.Range("I22").Value = [double]"0.126087115052153"
$ws.Range("J22").Value = [double]"0.126087115052153"
$ws.Range("M22").Value = [double]"6.305846"
$ws.Range("N22").Value = [double]"18.917538"
$ws.Range("O22").Value = [double]"0.01356150511917599"
$ws.Range("P22").Value = [double]"0.01356150511917599"
$ws.Range("Q22").Value = [double]"51.77443024412133"
$ws.Range("R22").Value = [double]"465.969872197092"
$ws.Range("S22").Value = [double]"0.001709931056241905"
$ws.Range("T22").Value = [double]"0.001709931056241905"
$ws.Range("G23").Value = [double]"8.210544666666665"
$ws.Range("H23").Value = [double]"24.631634"
$ws.Range("I23").Value = [double]"0.126087115052153"
$ws.Range("J23").Value = [double]"0.126087115052153"
$ws.Range("O23").Value = [double]"0.392557056479861"
$ws.Range("P23").Value = [double]"0.3925570564798609"
$ws.Range("Q23").Value = [double]"1498.684530879644"
$ws.Range("R23").Value = [double]"13488.1607779168"
$ws.Range("S23").Value = [double]"0.04949638674491076"
$ws.Range("T23").Value = [double]"0.04949638674491076"
$ws.Range("G24").Value = [double]"8.210544666666665"
$ws.Range("H24").Value = [double]"24.631634"
$ws.Range("I24").Value = [double]"0.126087115052153"
$ws.Range("J24").Value = [double]"0.126087115052153"
$ws.Range("M24").Value = [double]"127.396393"
$ws.Range("N24").Value = [double]"382.189179"
$ws.Range("O24").Value = [double]"0.2739817680029065"
$ws.Range("P24").Value = [double]"0.2739817680029065"
$ws.Range("Q24").Value = [double]"1045.993775098721"
$ws.Range("R24").Value = [double]"9413.943975888486"
$ws.Range("S24").Value = [double]"0.03454557070437476"
$ws.Range("T24").Value = [double]"0.03454557070437477"
$ws.Range("G25").Value = [double]"8.210544666666665"
$ws.Range("H25").Value = [double]"24.631634"
$ws.Range("I25").Value = [double]"0.126087115052153"
$ws.Range("J25").Value = [double]"0.126087115052153"
$ws.Range("M25").Value = [double]"19.42400133333333"
$ws.Range("N25").Value = [double]"58.272004"
$ws.Range("O25").Value = [double]"0.04177372766745037"
$ws.Range("P25").Value = [double]"0.04177372766745036"
$ws.Range("Q25").Value = [double]"159.4816305527262"
$ws.Range("R25").Value = [double]"1435.334674974536"
$ws.Range("S25").Value = [double]"0.005267128806563122"
$ws.Range("T25").Value = [double]"0.005267128806563122"
$ws.Range("G26").Value = [double]"8.210544666666665"
$ws.Range("H26").Value = [double]"24.631634"
$ws.Range("I26").Value = [double]"0.126087115052153"
$ws.Range("J26").Value = [double]"0.126087115052153"
$ws.Range("M26").Value = [double]"129.3233566666667"
$ws.Range("N26").Value = [double]"387.97007"
$ws.Range("O26").Value = [double]"0.2781259427306063"
$ws.Range("P26").Value = [double]"0.2781259427306062"
$ws.Range("Q26").Value = [double]"1061.815196354931"
$ws.Range("R26").Value = [double]"9556.33676719438"
$ws.Range("S26").Value = [double]"0.03506809774006247"
$ws.Range("T26").Value = [double]"0.03506809774006247"
